$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.228.98'
$ws.Range("E2").Value = '  +8.23%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.593.62'
$ws.Range("E3").Value = '  +8.06%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.009'
$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9947'
$ws.Range("E5").Value = '  +4.00%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '301.40'
$ws.Range("E6").Value = '  +8.40%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3629'
$ws.Range("E7").Value = '  +0.22%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3363'
$ws.Range("E8").Value = '  +9.53%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '41.56'
$ws.Range("E9").Value = '  +4.75%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.115'
$ws.Range("E10").Value = '  +4.56%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06972'
$ws.Range("E11").Value = '  +4.61%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.008'
$ws.Range("E12").Value = '  +0.57%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.850'
$ws.Range("E13").Value = '  +5.71%  '

$ws.Range("B14").Value = 'Solana'
$ws.Range("C14").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.40'
$ws.Range("E14").Value = '  +6.78%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.566'
$ws.Range("E15").Value = '  +5.91%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9962'
$ws.Range("E16").Value = '  +4.21%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.593.71'
$ws.Range("E17").Value = '  +7.97%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001064'
$ws.Range("E18").Value = '  +3.51%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06601'
$ws.Range("E19").Value = '  +11.43%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '76.40'
$ws.Range("E20").Value = '  +10.57%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.943'
$ws.Range("E21").Value = '  +7.78%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.76'
$ws.Range("E22").Value = '  +8.48%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.65'
$ws.Range("E23").Value = '  +4.21%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.302.26'
$ws.Range("E24").Value = '  +8.50%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.397'
$ws.Range("E25").Value = '  +6.54%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.507'
$ws.Range("E26").Value = '  +17.28%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '148.27'
$ws.Range("E27").Value = '  +3.42%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.16'
$ws.Range("E28").Value = '  +11.51%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.769.06'
$ws.Range("E29").Value = '  +8.07%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '121.70'
$ws.Range("E30").Value = '  +6.88%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.053'
$ws.Range("E31").Value = '  +3.33%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.941'
$ws.Range("E32").Value = '  +19.16%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9229'
$ws.Range("E33").Value = '  +13.89%  '

$ws.Range("B34").Value = 'Stellar'
$ws.Range("C34").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08174'
$ws.Range("E34").Value = '  +2.18%  '

$ws.Range("B35").Value = 'WEMIXTOKEN'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.640'
$ws.Range("E35").Value = '  +8.24%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '11.70'
$ws.Range("E36").Value = '  +12.57%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.119'
$ws.Range("E37").Value = '  +7.96%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.237'
$ws.Range("E38").Value = '  +1.62%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.371'
$ws.Range("E39").Value = '  +12.65%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06006'
$ws.Range("E40").Value = '  +3.58%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02185'
$ws.Range("E41").Value = '  +6.14%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1987'
$ws.Range("E42").Value = '  +5.62%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9946'
$ws.Range("E43").Value = '  +3.94%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5801'
$ws.Range("E44").Value = '  +9.64%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.794'
$ws.Range("E45").Value = '  +7.70%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.92'
$ws.Range("E46").Value = '  +5.41%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5588'
$ws.Range("E47").Value = '  +7.27%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.10'
$ws.Range("E48").Value = '  +5.62%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.943'
$ws.Range("E49").Value = '  +6.96%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06754'
$ws.Range("E50").Value = '  +4.25%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.59'
$ws.Range("E51").Value = '  +7.55%  '
